$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B / C / D: new/updated values (text + numbers + formulas) ---
# NOTE: cell writes are ordered to reproduce the original author's shared-string
# insertion order (so new shared-string indices line up with the target file).

# Row 5: Current/voltage regulated answer
$ws.Range("B5").Value = "Voltage Regulated"

# Row 6: Maximum Output Voltage = 0.001 * 510 (formula)
$ws.Range("B6").Formula = "=0.001*510"

# Row 7: Maximum Output Current = 0.001 (plain number)
$ws.Range("B7").Value = 0.001

# Row 9: Pulse Duration = 0.00025 (plain number)
$ws.Range("B9").Value = 0.00025

# Row 13: Net Charge note
$ws.Range("B13").Value = "Current / pulse width"

# Row 14: Leakage Current note + formula in C14
$ws.Range("B14").Value = "Peak voltage / resistance"
$ws.Range("C14").Formula = "=((0.862-0.532)/510)*10^9"

# Row 15: Net DC Current note + formula in C15
$ws.Range("B15").Value = "Max volt peak/510- min volt/510"
$ws.Range("C15").Formula = "=((0.532-0.528)/510)*10^6"

# Row 16: Maximum Phase Charge formula
$ws.Range("B16").Formula = "=(0.532/510)*(250/(10^6))"

# Row 18: Maximum Phase Power formula
$ws.Range("B18").Formula = "=0.532*(0.001)"

# Row 19: Maximum Phase Power Density note
$ws.Range("B19").Value = "Use cell B17 and divide by elerode area"

# Row 20: Pulse Delivery Mode - fix "Burs" -> "Burst"
$ws.Range("B20").Value = "Burst"

# Row 14 marker (D14)
$ws.Range("D14").Value = "240hz, 250 uS, 1 mA"

# Row 13 marker (D13)
$ws.Range("D13").Value = "Default"

# Row 10: Frequency range
$ws.Range("B10").Value = "30-240"

# Row 8: multiphasic waveform answer
$ws.Range("B8").Value = "Biphasic symmetrical"

# --- Apply Wrap Text formatting to the whole column B (matches all B cells getting style s="1") ---
$ws.Columns("B").WrapText = $true

# --- Row heights: row 5 shrinks from 45 to 30 ---
$ws.Rows(5).RowHeight = 30

# --- Column C width ---
$ws.Columns("C").ColumnWidth = 11.1

# --- Selection / active cell to match final view ---
$ws.Range("B22").Select() | Out-Null
